$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A7").Value = "KIRAN KUMAR"
$ws.Range("B7").Value = "OS"
$ws.Range("C7").Value = "Ftth OS_01.12.2025.xlsx"
$ws.Range("D7").Value = "2025-12-02 12:19"
$ws.Range("E7").Value = "2025-12"
